$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.238.31"
$ws.Range("E2").Value = "  -3.52%  "
$ws.Range("D3").Value = "2.981.64"
$ws.Range("E3").Value = "  -3.11%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.56"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.16"
$ws.Range("E6").Value = "  -5.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "2.978.40"
$ws.Range("E8").Value = "  -3.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.145"
$ws.Range("E10").Value = "  -5.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.94"
$ws.Range("E11").Value = "  -7.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.443"
$ws.Range("E12").Value = "  -2.49%  "
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.67"
$ws.Range("E14").Value = "  -2.74%  "
$ws.Range("D15").Value = "3.462.99"
$ws.Range("E15").Value = "  -3.32%  "
$ws.Range("D16").Value = "61.232.68"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("E17").Value = "  -2.54%  "
$ws.Range("D18").Value = "2.980.02"
$ws.Range("E18").Value = "  -3.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.56"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "469.76"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.01"
$ws.Range("E21").Value = "  -2.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.664"
$ws.Range("E22").Value = "  -4.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.94"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.56"
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.89"
$ws.Range("E25").Value = "  -2.68%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.69"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.59"
$ws.Range("E28").Value = "  -5.13%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.88"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.36"
$ws.Range("E31").Value = "  -3.05%  "
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.45"
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.28"
$ws.Range("E34").Value = "  -2.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "54.81"
$ws.Range("E35").Value = "  -4.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.83"
$ws.Range("E36").Value = "  -2.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "444.62"
$ws.Range("E37").Value = "  -9.10%  "
$ws.Range("D38").Value = "3.134.19"
$ws.Range("E38").Value = "  -3.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0784"
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0376"
$ws.Range("E40").Value = "  -6.13%  "
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.03"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").Value = "  -11.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.63"
$ws.Range("E45").Value = "  +2.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.240"
$ws.Range("E46").Value = "  -5.12%  "
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.93"
$ws.Range("E48").Value = "  -4.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "114.86"
$ws.Range("E49").Value = "  -7.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.30"
$ws.Range("E50").Value = "  +8.96%  "
$ws.Range("D51").Value = "0.0₃0480"
$ws.Range("E51").Value = "  -9.19%  "
